# Swap the contents of column D (codeforiati:category-name / category values)
# and column G (codeforiati:group-code / group values) for every row, including
# the header row. Columns A, B, C, E, F are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$colD = 4
$colG = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, $colD)
    $gCell = $ws.Cells.Item($r, $colG)

    $dVal = $dCell.Value()
    $gVal = $gCell.Value()

    $dCell.Value2 = $gVal
    $gCell.Value2 = $dVal
}
